# Generate Report for Handoff
# A fresh localization-status report run: the zh-cn / de-de files moved from
# "In Translation" to "Ready for handoff", and the HO Xliff generation /
# handoff timestamps were refreshed to the new run's datetime.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# zh-cn / de-de status columns (E, F) and the "Latest HO Xliff Generate Date"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-04 23:03:51"

# --- zh-cn sheet ------------------------------------------------------------
# Status column (C) and Latest Handoff Datetime (H)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 23:03:47"

# --- de-de sheet ------------------------------------------------------------
# Status column (C) and Latest Handoff Datetime (H)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-04 23:03:51"

# The new status text ("Ready for handoff") is longer than the old one
# ("In Translation"), so the status columns widen to fit it.
$wsOverview.Columns("E:E").ColumnWidth = 16.333333333333336
$wsOverview.Columns("F:F").ColumnWidth = 16.333333333333336
$wsZhCn.Columns("C:C").ColumnWidth = 16.333333333333336
$wsDeDe.Columns("C:C").ColumnWidth = 16.333333333333336
